$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure price/volume cells stay text so formatted strings
# (thousand-dot separators, trailing zeros, padded %) survive
# instead of being auto-coerced to numbers by COM assignment.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "69.583.18"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -0.39%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.501.40"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -0.59%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "574.85"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -0.50%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "166.60"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -1.03%  "
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -1.54%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.500.07"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -0.59%  "
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -0.42%  "
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -0.02%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.356"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +2.96%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.958.71"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -0.37%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "69.496.14"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -0.49%  "
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +0.86%  "
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -1.47%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.505.52"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -0.28%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.18"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -1.11%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.44"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -4.55%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "348.15"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -0.55%  "
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -0.76%  "
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -0.64%  "
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -0.17%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "70.75"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +2.80%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.93"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -1.56%  "
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -2.54%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.629.96"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -0.52%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.992"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -0.65%  "
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -1.93%  "
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -0.70%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "457.72"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -2.03%  "
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -4.73%  "
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -1.59%  "
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +0.03%  "
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -3.31%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "156.58"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +2.28%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "19.04"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +0.24%  "
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -0.41%  "
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +0.01%  "
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -0.98%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "4.68"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -2.10%  "
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -0.52%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "38.08"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -0.26%  "
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -5.03%  "
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -6.53%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "141.05"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -1.61%  "
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -0.51%  "
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -2.52%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0732"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -0.56%  "
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -0.73%  "
